$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "Förändrad" (last changed) date, stored as a serial
# date number. Every data row (2 through 110) was bumped by one day,
# from 45189 (2023-09-20) to 45190 (2023-09-21).
$ws.Range("C2:C110").Value = 45190
